## added test for sortMyPolicyPage
## Appends a new "sortPolicyList" data block to the DashboardPageData sheet,
## mirroring the existing "sortQuoteList" block (rows 55-58) three rows lower
## (rows 62-65), leaving rows 59-61 blank just like the gap before row 55.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DashboardPageData")
$ws.Activate()

# Copy the existing "sortQuoteList" block (title + header + 2 data rows)
# down to rows 62-65, preserving styles/borders/fills exactly.
$ws.Range("A55:B55").Copy($ws.Range("A62:B62"))
$ws.Range("A56:D56").Copy($ws.Range("A63:D63"))
$ws.Range("A57:D57").Copy($ws.Range("A64:D64"))
$ws.Range("A58:D58").Copy($ws.Range("A65:D65"))

# Rename the new block's title from "sortQuoteList" to "sortPolicyList".
$ws.Range("A62").Value = "sortPolicyList"

# Match the author's final selection, as left after adding the new rows.
$null = $ws.Range("A62").Select()
